$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 239
$ws.Range("I19").Value = 168
$ws.Range("K19").Value = 168
$ws.Range("M19").Value = 7
$ws.Range("H100").Value = 3175.75
$ws.Range("I100").Value = 2000
$ws.Range("K100").Value = 2000
$ws.Range("M100").Value = -1459
$ws.Range("H129").Value = 239030.17
$ws.Range("J129").Value = 271293.7
$ws.Range("L129").Value = 813881.1000000001
$ws.Range("N129").Value = -823881.1000000001
$ws.Range("H137").Value = 97714.59
$ws.Range("I137").Value = 123199.15
$ws.Range("J137").Value = 4271.222
$ws.Range("K137").Value = 369597.45
$ws.Range("L137").Value = 12813.666
$ws.Range("M137").Value = -367047.45
$ws.Range("N137").Value = -17913.666
$ws.Range("H138").Value = 3899.6438
$ws.Range("J138").Value = 3713.0896
$ws.Range("L138").Value = 11139.2688
$ws.Range("N138").Value = -21419.2688

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2117.0908
$ws.Range("I2").Value = 2181.3333
$ws.Range("J2").Value = 2040
$ws.Range("K2").Value = 2181.3333
$ws.Range("L2").Value = 2040
$ws.Range("M2").Value = -2068.3333
$ws.Range("N2").Value = -2266
$ws.Range("H116").Value = 2117.0908
$ws.Range("I116").Value = 2181.3333
$ws.Range("J116").Value = 2040
$ws.Range("K116").Value = 2181.3333
$ws.Range("L116").Value = 2040
$ws.Range("M116").Value = 112.6667000000002
$ws.Range("N116").Value = -6628

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2117.0908
$ws.Range("I3").Value = 2181.3333
$ws.Range("J3").Value = 2040
$ws.Range("K3").Value = 2181.3333
$ws.Range("L3").Value = 2040
$ws.Range("M3").Value = -2067.3333
$ws.Range("N3").Value = -2268
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("H26").Value = 12290.333
$ws.Range("I26").Value = 3435.5
$ws.Range("K26").Value = 3435.5
$ws.Range("M26").Value = -3143.5
$ws.Range("H36").Value = 768.5
$ws.Range("I36").Value = 768.5
$ws.Range("K36").Value = 768.5
$ws.Range("M36").Value = -234.5
$ws.Range("H54").Value = 9651.777
$ws.Range("I54").Value = 3773.2
$ws.Range("K54").Value = 3773.2
$ws.Range("M54").Value = -3289.2
$ws.Range("H82").Value = 27511.875
$ws.Range("I82").Value = 7817.5
$ws.Range("J82").Value = 47206.25
$ws.Range("K82").Value = 7817.5
$ws.Range("L82").Value = 47206.25
$ws.Range("M82").Value = -7434.5
$ws.Range("N82").Value = -47972.25
$ws.Range("H85").Value = 27511.875
$ws.Range("I85").Value = 7817.5
$ws.Range("J85").Value = 47206.25
$ws.Range("K85").Value = 7817.5
$ws.Range("L85").Value = 47206.25
$ws.Range("M85").Value = -6491.5
$ws.Range("N85").Value = -49858.25
$ws.Range("H96").Value = 26854.2
$ws.Range("I96").Value = 1300
$ws.Range("J96").Value = 33242.75
$ws.Range("K96").Value = 1300
$ws.Range("L96").Value = 33242.75
$ws.Range("M96").Value = 1446
$ws.Range("N96").Value = -38734.75
$ws.Range("H97").Value = 14564.728
$ws.Range("J97").Value = 25000
$ws.Range("L97").Value = 25000
$ws.Range("N97").Value = -26982
$ws.Range("H134").Value = 3286.9636
$ws.Range("I134").Value = 2985
$ws.Range("J134").Value = 5061
$ws.Range("K134").Value = 8955
$ws.Range("L134").Value = 15183
$ws.Range("M134").Value = -6420
$ws.Range("N134").Value = -20253
$ws.Range("M24").ClearContents()

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H43").Value = 27661.666
$ws.Range("J43").Value = 27661.666
$ws.Range("L43").Value = 27661.666
$ws.Range("N43").Value = -28029.666
$ws.Range("H94").Value = 4512.4165
$ws.Range("I94").Value = 3055.8333
$ws.Range("J94").Value = 5969
$ws.Range("K94").Value = 3055.8333
$ws.Range("L94").Value = 5969
$ws.Range("M94").Value = -2604.8333
$ws.Range("N94").Value = -6871
$ws.Range("H99").Value = 3755.1853
$ws.Range("I99").Value = 2874.5
$ws.Range("K99").Value = 2874.5
$ws.Range("M99").Value = -1376.5
$ws.Range("H101").Value = 27661.666
$ws.Range("J101").Value = 27661.666
$ws.Range("L101").Value = 27661.666
$ws.Range("N101").Value = -34151.666
$ws.Range("H107").Value = 2096.9443
$ws.Range("I107").Value = 1374
$ws.Range("J107").Value = 2675.3
$ws.Range("K107").Value = 1374
$ws.Range("L107").Value = 2675.3
$ws.Range("M107").Value = 546
$ws.Range("N107").Value = -6515.3
$ws.Range("H108").Value = 30621
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("H124").Value = 10318.434
$ws.Range("I124").Value = 8162.577
$ws.Range("J124").Value = 24331.5
$ws.Range("K124").Value = 8162.577
$ws.Range("L124").Value = 24331.5
$ws.Range("M124").Value = -5707.577
$ws.Range("N124").Value = -29241.5
$ws.Range("H126").Value = 3755.1853
$ws.Range("I126").Value = 2874.5
$ws.Range("K126").Value = 8623.5
$ws.Range("M126").Value = -6153.5
$ws.Range("N108").ClearContents()

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 90
$ws.Range("I98").Value = 90
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 270
$ws.Range("L98").Value = 0
$ws.Range("N98").Value = 1228
$ws.Range("H121").Value = 787.3514
$ws.Range("I121").Value = 174.91667
$ws.Range("J121").Value = 1081.32
$ws.Range("K121").Value = 524.75001
$ws.Range("L121").Value = 3243.96
$ws.Range("M121").Value = 785.24999
$ws.Range("N121").Value = -5863.96
$ws.Range("M98").ClearContents()

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4127
$ws.Range("I80").Value = 3883.3333
$ws.Range("J80").Value = 4289.4443
$ws.Range("K80").Value = 3883.3333
$ws.Range("L80").Value = 4289.4443
$ws.Range("M80").Value = -2885.3333
$ws.Range("N80").Value = -6285.4443
$ws.Range("H83").Value = 4127
$ws.Range("I83").Value = 3883.3333
$ws.Range("J83").Value = 4289.4443
$ws.Range("K83").Value = 19416.6665
$ws.Range("L83").Value = 21447.2215
$ws.Range("M83").Value = -14424.6665
$ws.Range("N83").Value = -31431.2215
$ws.Range("H95").Value = 25000
$ws.Range("J95").Value = 25000
$ws.Range("L95").Value = 25000
$ws.Range("N95").Value = -30492
$ws.Range("H99").Value = 8000
$ws.Range("I99").Value = 8000
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 8000
$ws.Range("L99").Value = 0
$ws.Range("N99").Value = -5754
$ws.Range("H102").Value = 4977.1816
$ws.Range("I102").Value = 4859.4443
$ws.Range("K102").Value = 4859.4443
$ws.Range("M102").Value = -3237.4443
$ws.Range("H113").Value = 11688.875
$ws.Range("I113").Value = 17982.2
$ws.Range("K113").Value = 17982.2
$ws.Range("M113").Value = -15812.2
$ws.Range("M99").ClearContents()

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 275.9
$ws.Range("I55").Value = 270.6
$ws.Range("J55").Value = 291.8
$ws.Range("K55").Value = 270.6
$ws.Range("L55").Value = 291.8
$ws.Range("M55").Value = -97.60000000000002
$ws.Range("N55").Value = -637.8
$ws.Range("H122").Value = 1816079.5
$ws.Range("J122").Value = 5100
$ws.Range("L122").Value = 15300
$ws.Range("N122").Value = -20200

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H95").Value = 17800
$ws.Range("J95").Value = 17800
$ws.Range("L95").Value = 17800
$ws.Range("N95").Value = -23292
$ws.Range("H136").Value = 47624024
$ws.Range("I136").Value = 100002850
$ws.Range("J136").Value = 6909.5454
$ws.Range("K136").Value = 300008550
$ws.Range("L136").Value = 20728.6362
$ws.Range("M136").Value = -300006000
$ws.Range("N136").Value = -25828.6362
